$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B/C column text updates (swap Kaspa/Aptos rows 30 and 31) ---
$ws.Range("B30").Value = 'Aptos'
$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'

# --- D column (Price) updates: force text storage to preserve exact formatting ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.897.49'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.615.55'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.74'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.10'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.82'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.384'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.30'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.086.18'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000183'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.623.34'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.591.30'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.08'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.65'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.58'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '343.36'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.23'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.75'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000110'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '596.74'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.22'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.58'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.93'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.160'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.74'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.60'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.37'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.80'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '154.58'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.87'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '41.55'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.47'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '156.02'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.93'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.06'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0593'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.102'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.628'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0248'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.94'

# --- E column (Volume 1h) updates: force text storage to preserve exact formatting ---
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.78%  '
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.35%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -3.46%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.90%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -4.91%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.10%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.46%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.24%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.07%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.37%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -6.06%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.99%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.59%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.70%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.96%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.79%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.69%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.12%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.50%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +6.78%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.18%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.94%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.30%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.25%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.90%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.63%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.07%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.01%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.04%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.56%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.86%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.15%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.46%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.10%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.59%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +8.52%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.10%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.55%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.31%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.07%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.23%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.36%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.91%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.85%  '
